$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 317-328 (Fruta / Comercializadora del Agro de Limari - Frutilla)
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo, G Producto ID,
#          H Producto, I Categoria ID, J Categoria, K Variedad, L Calidad, M Volumen,
#          N Precio minimo, O Precio maximo, P Precio promedio ponderado,
#          Q Unidad de comercializacion, R Origen, S Precio $/Kg, T Kg/unidad
$rows = @(
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44595, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 300, 11500, 12000, 11750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1679, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44595, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 400, 9500, 10000, 9750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1393, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44595, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 300, 7500, 8000, 7750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1107, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44335, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 240, 19500, 20000, 19750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 2821, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44335, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 240, 16500, 17000, 16750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 2393, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44335, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 300, 11500, 12000, 11750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1679, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44552, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 400, 12500, 13000, 12750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1821, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44552, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 340, 10500, 11000, 10750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1536, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44552, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 260, 8500, 9000, 8750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1250, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44160, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 300, 14500, 15000, 14750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 2107, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44160, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 300, 12500, 13000, 12750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1821, 7),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44160, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 240, 10500, 11000, 10750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1536, 7),
)

$startRow = 317
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $col = [char](65 + $c)
        $ws.Range("$col$r").Value = $data[$c]
    }
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
